# Generate Report for Handoff
# The 9c9bafea-51f4-4c8b-a552-1c000fb1ea13.md file has been re-handed-off;
# update its status/dates on the Overview, zh-cn and de-de sheets and
# record the "stale handback" error detail message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce1d58183d77909a207c2cea464ad0cd29572212/e2e/9c9bafea-51f4-4c8b-a552-1c000fb1ea13.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49391883f2c279cd44a047fdd7d9be6c9b474309/e2e/9c9bafea-51f4-4c8b-a552-1c000fb1ea13.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 14:52:48"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-04 14:52:43"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-04 14:52:48"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1
